$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Reference cells used as formatting "donors" so the new sheet reuses the
# workbook's existing styles instead of growing the style table.
$headerDonor = $wb.Worksheets.Item(1).Range("A1")   # bold/border/center header style
$plainDonor  = $wb.Worksheets.Item(1).Range("B2")   # default/no style

function Set-TextCell($cell, [string]$text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (even if it looks numeric, or is empty) instead of inferring a number.
    $cell.Formula = "'" + $text
    # The apostrophe entry also stamps a transient "quote prefix" style on
    # the cell; repaint it with the donor's plain style so the saved file
    # doesn't carry that extra formatting.
    $plainDonor.Copy()
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Header row ---
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$headerDonor.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data row ---
Set-TextCell $newSheet.Range("A2") "4485"
$newSheet.Range("B2").Value = 11
Set-TextCell $newSheet.Range("C2") "0"
Set-TextCell $newSheet.Range("D2") "0"
Set-TextCell $newSheet.Range("E2") ""
Set-TextCell $newSheet.Range("F2") "NO"

$excel.CutCopyMode = 0
